$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.676.72"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.099.01"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5158"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4383"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09180"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.35%  "
$ws.Range("D13").Value = "2.092.48"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.762"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.189"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001152"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06676"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.209"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "29.737.30"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "2.337.07"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.498"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.128"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.98%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.668"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1050"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.205"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.304"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02577"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06679"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6986"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.330"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6824"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.309"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.615"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000353"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.206"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.217"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.93%  "
